$wb = $excel.ActiveWorkbook

# --- Todo sheet: append a new to-do item ---
$todo = $wb.Worksheets.Item("Todo")
$todo.Range("A5").Value = "I guess I added wrong events to 372 or 373"

# --- GotoPoints sheet: append the new Cavetown goto points ---
$goto = $wb.Worksheets.Item("GotoPoints")
$goto.Range("A6").Value  = "79: Warenhändler / Good merchant (Cavetown)"
$goto.Range("A7").Value  = "80: Schmied / Blacksmith (Cavetown)"
$goto.Range("A8").Value  = "81: Cavetown Büro / Cavetown Office (Cavetown)"
$goto.Range("A9").Value  = "82: Vielauge-Schloss / Manyeyes' Castle (Cavetown)"
$goto.Range("A10").Value = "83: Gasthaus / Tavern (Cavetown)"
$goto.Range("A11").Value = "84: Badehaus / Bathhouse (Cavetown)"
$goto.Range("A12").Value = "85: Flosshändler / Raft Dealer (Cavetown)"
$goto.Range("A13").Value = "86: Gasthaus / Tavern (Cavetown) -- Second door"
$goto.Range("A14").Value = "87: Stadthaus 1 / Townhouse 1 (Cavetown)"
$goto.Range("A15").Value = "88: Stadthaus 2 / Townhouse 2 (Cavetown)"
$goto.Range("A16").Value = "89: Stadthaus 3 / Townhouse 3 (Cavetown)"

# --- Selections / active sheet bookkeeping ---
# Todo sheet ends up with its selection resting on B5 (no longer the active tab)
$todo.Activate()
$null = $todo.Range("B5").Select()

# GotoPoints becomes the active/selected tab, with the selection on C15
$goto.Activate()
$null = $goto.Range("C15").Select()
